$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Torta (row 2): reorder raw materials list, Vainilla moved to end
$ws.Range("C2").Value = "1.0-Huevos (unidad),3.0-Leche (litros),2.0-Harina  (kg),1.0-Vainilla (ml),"

# Disponible flag switched from 0 to 1 (all rows shared the same underlying
# string value in the original workbook, so all of them flip together).
# Force text format so the value is stored as a string, matching the source.
$ws.Range("E2:E7").NumberFormat = "@"
$ws.Range("E2:E7").Value = "1"

# Queque (row 4): reorder raw materials list, Vainilla moved to end
$ws.Range("C4").Value = "2.0-Huevos (unidad),5.0-Harina  (kg),1.0-Vainilla (ml),"

# Pie de Limon (row 6): reorder raw materials list, Crema moved to front, Limon moved to end
$ws.Range("C6").Value = "5.0-Crema (litros),5.0-Huevos (unidad),4.0-Harina  (kg),2.0-Limon (unidad),"

# Cupcake (row 7): reorder raw materials list, Vainilla moved to end
$ws.Range("C7").Value = "2.0-Huevos (unidad),0.2-Leche (litros),0.3-Harina  (kg),0.1-Vainilla (ml),"
